$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 631.5333000000001
$ws.Range("J17").Value = 631.5333000000001
$ws.Range("L17").Value = 1894.5999
$ws.Range("N17").Value = -2230.5999
$ws.Range("H51").Value = 6877.9
$ws.Range("J51").Value = 7197.6665
$ws.Range("L51").Value = 7197.6665
$ws.Range("N51").Value = -8165.6665
$ws.Range("H98").Value = 783.9375
$ws.Range("I98").Value = 756.26666
$ws.Range("J98").Value = 1199
$ws.Range("K98").Value = 756.26666
$ws.Range("L98").Value = 1199
$ws.Range("M98").Value = 741.73334
$ws.Range("N98").Value = -4195
$ws.Range("H111").Value = 4283.143
$ws.Range("I111").Value = 5720.5
$ws.Range("J111").Value = 2366.6667
$ws.Range("K111").Value = 17161.5
$ws.Range("L111").Value = 7100.000100000001
$ws.Range("M111").Value = -14094.5
$ws.Range("N111").Value = -13234.0001
$ws.Range("H116").Value = 1988.4445
$ws.Range("J116").Value = 2401.5
$ws.Range("L116").Value = 2401.5
$ws.Range("N116").Value = -9285.5
$ws.Range("H122").Value = 783.9375
$ws.Range("I122").Value = 756.26666
$ws.Range("J122").Value = 1199
$ws.Range("K122").Value = 2268.79998
$ws.Range("L122").Value = 3597
$ws.Range("M122").Value = 181.2000200000002
$ws.Range("N122").Value = -8497
$ws.Range("H125").Value = 1353
$ws.Range("I125").Value = 906
$ws.Range("J125").Value = 1800
$ws.Range("K125").Value = 8154
$ws.Range("L125").Value = 16200
$ws.Range("M125").Value = -5694
$ws.Range("N125").Value = -21120
$ws.Range("H127").Value = 1035.2
$ws.Range("I127").Value = 554.8889
$ws.Range("J127").Value = 1755.6666
$ws.Range("K127").Value = 1664.6667
$ws.Range("L127").Value = 5266.9998
$ws.Range("M127").Value = 3295.3333
$ws.Range("N127").Value = -15186.9998

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3338
$ws.Range("I2").Value = 4326.6665
$ws.Range("J2").Value = 2914.2856
$ws.Range("K2").Value = 4326.6665
$ws.Range("L2").Value = 2914.2856
$ws.Range("M2").Value = -4213.6665
$ws.Range("N2").Value = -3140.2856
$ws.Range("H45").Value = 2300
$ws.Range("I45").Value = 2880
$ws.Range("J45").Value = 1816.6666
$ws.Range("K45").Value = 2880
$ws.Range("L45").Value = 1816.6666
$ws.Range("M45").Value = -2503
$ws.Range("N45").Value = -2570.6666
$ws.Range("H63").Value = 3675.625
$ws.Range("I63").Value = 3317.5
$ws.Range("K63").Value = 3317.5
$ws.Range("M63").Value = -2631.5
$ws.Range("H66").Value = 3675.625
$ws.Range("I66").Value = 3317.5
$ws.Range("K66").Value = 16587.5
$ws.Range("M66").Value = -13155.5
$ws.Range("H116").Value = 3338
$ws.Range("I116").Value = 4326.6665
$ws.Range("J116").Value = 2914.2856
$ws.Range("K116").Value = 4326.6665
$ws.Range("L116").Value = 2914.2856
$ws.Range("M116").Value = -2032.6665
$ws.Range("N116").Value = -7502.2856
$ws.Range("H122").Value = 4117194.8
$ws.Range("I122").Value = 1964.375
$ws.Range("K122").Value = 5893.125
$ws.Range("M122").Value = -3443.125

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3338
$ws.Range("I3").Value = 4326.6665
$ws.Range("J3").Value = 2914.2856
$ws.Range("K3").Value = 4326.6665
$ws.Range("L3").Value = 2914.2856
$ws.Range("M3").Value = -4212.6665
$ws.Range("N3").Value = -3142.2856
$ws.Range("H106").Value = 11866.667
$ws.Range("J106").Value = 11866.667
$ws.Range("L106").Value = 11866.667
$ws.Range("N106").Value = -14390.667
$ws.Range("H123").Value = 24998.572
$ws.Range("J123").Value = 24998.572
$ws.Range("L123").Value = 24998.572
$ws.Range("N123").Value = -34798.572

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 4906
$ws.Range("I94").Value = 15804
$ws.Range("J94").Value = 1273.3334
$ws.Range("K94").Value = 15804
$ws.Range("L94").Value = 1273.3334
$ws.Range("M94").Value = -15353
$ws.Range("N94").Value = -2175.3334

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 90910470
$ws.Range("I12").Value = 200002220
$ws.Range("J12").Value = 681.6667
$ws.Range("K12").Value = 600006660
$ws.Range("L12").Value = 2045.0001
$ws.Range("M12").Value = -600006487
$ws.Range("N12").Value = -2391.0001
$ws.Range("H113").Value = 600.6949
$ws.Range("H121").Value = 40089870
$ws.Range("I121").Value = 1003.75
$ws.Range("J121").Value = 48529636
$ws.Range("K121").Value = 3011.25
$ws.Range("L121").Value = 145588908
$ws.Range("M121").Value = -1701.25
$ws.Range("N121").Value = -145591528
$ws.Range("H131").Value = 1084.6666
$ws.Range("I131").Value = 467.14285
$ws.Range("J131").Value = 1166.2264
$ws.Range("K131").Value = 1401.42855
$ws.Range("L131").Value = 3498.6792
$ws.Range("M131").Value = 3638.57145
$ws.Range("N131").Value = -13578.6792
$ws.Range("H132").Value = 1090.4736
$ws.Range("I132").Value = 678.46155
$ws.Range("J132").Value = 1983.1666
$ws.Range("K132").Value = 6106.15395
$ws.Range("L132").Value = 17848.4994
$ws.Range("M132").Value = -3576.15395
$ws.Range("N132").Value = -22908.4994
$ws.Range("H138").Value = 2762.6924
$ws.Range("I138").Value = 1691.3334
$ws.Range("J138").Value = 4223.636
$ws.Range("K138").Value = 5074.0002
$ws.Range("L138").Value = 12670.908
$ws.Range("M138").Value = 65.9997999999996
$ws.Range("N138").Value = -22950.908

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 22675
$ws.Range("I29").Value = 700
$ws.Range("J29").Value = 30000
$ws.Range("K29").Value = 700
$ws.Range("L29").Value = 30000
$ws.Range("M29").Value = -410
$ws.Range("N29").Value = -30580
$ws.Range("H70").Value = 28844.879
$ws.Range("I70").Value = 42536.152
$ws.Range("J70").Value = 5113.3335
$ws.Range("K70").Value = 42536.152
$ws.Range("L70").Value = 5113.3335
$ws.Range("M70").Value = -42266.152
$ws.Range("N70").Value = -5653.3335
$ws.Range("H73").Value = 28844.879
$ws.Range("I73").Value = 42536.152
$ws.Range("J73").Value = 5113.3335
$ws.Range("K73").Value = 42536.152
$ws.Range("L73").Value = 5113.3335
$ws.Range("M73").Value = -41600.152
$ws.Range("N73").Value = -6985.3335
$ws.Range("H113").Value = 1570.6666
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1570.6666
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1570.6666
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -5910.6666
$ws.Range("H122").Value = 3134.8125
$ws.Range("I122").Value = 2810.4666
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 8431.399800000001
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -5981.399800000001
$ws.Range("N122").Value = -28900
$ws.Range("H123").Value = 27401
$ws.Range("J123").Value = 27401
$ws.Range("L123").Value = 27401
$ws.Range("N123").Value = -32301
$ws.Range("H132").Value = 102411.9
$ws.Range("I132").Value = 78633.84
$ws.Range("K132").Value = 235901.52
$ws.Range("M132").Value = -233371.52
$ws.Range("H141").Value = 12714.546
$ws.Range("J141").Value = 18042.857
$ws.Range("L141").Value = 18042.857
$ws.Range("N141").Value = -28402.857

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 53083.566
$ws.Range("I132").Value = 23612.875
$ws.Range("K132").Value = 70838.625
$ws.Range("M132").Value = -68308.625

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 465.8889
$ws.Range("I107").Value = 200
$ws.Range("J107").Value = 499.125
$ws.Range("K107").Value = 600
$ws.Range("L107").Value = 1497.375
$ws.Range("M107").Value = 1320
$ws.Range("N107").Value = -5337.375
$ws.Range("H122").Value = 2083.919
$ws.Range("I122").Value = 1617.4348
$ws.Range("J122").Value = 2850.2856
$ws.Range("K122").Value = 4852.3044
$ws.Range("L122").Value = 8550.856800000001
$ws.Range("M122").Value = -2402.3044
$ws.Range("N122").Value = -13450.8568
$ws.Range("H132").Value = 73058.46000000001
$ws.Range("I132").Value = 72668.14
$ws.Range("J132").Value = 73448.78999999999
$ws.Range("K132").Value = 218004.42
$ws.Range("L132").Value = 220346.37
$ws.Range("M132").Value = -215474.42
$ws.Range("N132").Value = -225406.37
